# Apply the "change outcome names to values and labels" edit to the
# outcome table: several outcome_name labels (column B) are simplified /
# consolidated, and the active selection moves to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcome")

# --- Column B (outcome_name) label updates ---------------------------------
# c11_youth_health_costs / c16_youth_health_costs: drop the age-group
# qualifier now that it is already conveyed by the "population" column.
$ws.Range("B16").Value = "Zorgkosten"
$ws.Range("B17").Value = "Zorgkosten"

# c11_youth_protection / c16_youth_protection: same consolidation.
$ws.Range("B18").Value = "Jeugdbescherming "
$ws.Range("B19").Value = "Jeugdbescherming "

# c11_living_space_pp / c16_living_space_pp: same consolidation.
$ws.Range("B43").Value = "Woonoppervlak per lid huishouden"
$ws.Range("B44").Value = "Woonoppervlak per lid huishouden"

# --- Row heights: let Excel re-fit the rows to the (updated) content -------
$ws.Range("A1:H44").EntireRow.AutoFit()
$ws.Range("A4:H4").RowHeight = 32
$ws.Range("A6:H6").RowHeight = 32
$ws.Range("A20:H20").RowHeight = 32

# --- Selection moves to B1 before the file is saved -------------------------
$ws.Activate() | Out-Null
$ws.Range("B1").Select() | Out-Null
